$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# fruit_name "Republic" -> "Shirt" (row 5)
$ws.Range("B5").Value = "Shirt"

# price in D2: 299 -> 350
$ws.Range("D2").Value = 350

# Give every data row (1-7) an explicit 16pt custom row height
$ws.Range("A1:E7").RowHeight = 16

# Move the active selection from B4 to C2
$ws.Range("C2").Select()
